$wb = $excel.ActiveWorkbook

# --- Sheet "450_curated babbar": add row 25 ---
$ws1 = $wb.Worksheets.Item("450_curated babbar")
$ws1.Cells.Item(25, 1).Value = "18 feb"
$ws1.Cells.Item(25, 3).Value = "sort matrix"

# --- Sheet "mistakes": add row 22 ---
$ws3 = $wb.Worksheets.Item("mistakes")
$ws3.Cells.Item(22, 2).Value = "segmentation fault"
$ws3.Cells.Item(22, 3).Value = "check return values are provided`nCheck out of bound values are called?"
$ws3.Rows.Item(22).RowHeight = 24.45

# --- Selection / active sheet state ---
$ws1.Range("D25").Select()
$ws3.Range("C22").Select()
$ws3.Activate()
